$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.198.38"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "2.373.63"
$ws.Range("E3").Value = "  +2.92%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "303.59"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").Value = "97.27"
$ws.Range("E6").Value = "  +1.30%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").Value = "0.502"
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("D11").Value = "0.0791"
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("E12").Value = "  +3.53%  "
$ws.Range("D13").Value = "18.61"
$ws.Range("E13").Value = "  -2.98%  "
$ws.Range("D14").Value = "6.79"
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("D15").Value = "2.740.71"
$ws.Range("E15").Value = "  +2.66%  "
$ws.Range("D16").Value = "2.332.21"
$ws.Range("E16").Value = "  +1.22%  "
$ws.Range("D17").Value = "0.802"
$ws.Range("E17").Value = "  +2.12%  "
$ws.Range("D18").Value = "43.180.11"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").Value = "12.37"
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("D20").Value = "6.31"
$ws.Range("E20").Value = "  +4.78%  "
$ws.Range("D21").Value = "0.0₃0890"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "68.29"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("D23").Value = "236.11"
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("E24").Value = "  -2.76%  "
$ws.Range("D25").Value = "2.44"
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "24.81"
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("E29").Value = "  +0.95%  "
$ws.Range("D30").Value = "31.54"
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "5.12"
$ws.Range("E31").Value = "  +2.71%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").Value = "0.0728"
$ws.Range("E33").Value = "  +3.77%  "
$ws.Range("D34").Value = "17.33"
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("E35").Value = "  +5.59%  "
$ws.Range("D36").Value = "4.39"
$ws.Range("E36").Value = "  -1.24%  "
$ws.Range("E37").Value = "  -1.10%  "
$ws.Range("E38").Value = "  +1.38%  "
$ws.Range("D39").Value = "22.67"
$ws.Range("E39").Value = "  +12.72%  "
$ws.Range("D40").Value = "2.80"
$ws.Range("E40").Value = "  +3.44%  "
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("D42").Value = "1.946.48"
$ws.Range("E42").Value = "  -1.23%  "
$ws.Range("D43").Value = "102.71"
$ws.Range("E43").Value = "  -38.01%  "
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("E45").Value = "  +4.31%  "
$ws.Range("D46").Value = "9.44"
$ws.Range("E46").Value = "  -9.91%  "
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").Value = "2.596.47"
$ws.Range("E48").Value = "  +2.50%  "
$ws.Range("D49").Value = "53.14"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("E50").Value = "  +1.81%  "
$ws.Range("D51").Value = "72.37"
$ws.Range("E51").Value = "  +1.24%  "
